$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewTaxReturn")

# --- Cell edits on the "NewTaxReturn" sheet (Month -> CRN 'Y' flags + Year fixes) ---

# New "Y" (CRN) flags added in column H for several rows
$ws.Range("H21").Value = "Y"
$ws.Range("H24").Value = "Y"
$ws.Range("H25").Value = "Y"
$ws.Range("H26").Value = "Y"
$ws.Range("H33").Value = "Y"
$ws.Range("H36").Value = "Y"
$ws.Range("H37").Value = "Y"
$ws.Range("H38").Value = "Y"
$ws.Range("H45").Value = "Y"
$ws.Range("H48").Value = "Y"
$ws.Range("H49").Value = "Y"
$ws.Range("H50").Value = "Y"

# Year (column F) corrections
$ws.Range("F27").Value = "2024"
$ws.Range("F39").Value = "2023"
$ws.Range("F51").Value = "2022"

# --- View/selection state ---
# Personal_IND keeps its previous selection (E11) but loses the active tab
$ws2 = $wb.Worksheets.Item("Personal_IND")
$ws2.Select() | Out-Null
$ws2.Range("E11").Select() | Out-Null

# NewTaxReturn becomes the active tab, with G54 selected
$ws.Select() | Out-Null
$ws.Range("G54").Select() | Out-Null
